$d = $word.ActiveDocument

# The paragraph that currently ends with "...avoid mistakes." and carries
# the "_GoBack" bookmark at its end.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*avoid mistakes.*") {
        $target = $p
        break
    }
}

# Insert a brand-new (empty) list paragraph right after it; it inherits the
# same pPr/list numbering/rPr as the paragraph it follows.
$null = $target.Range.InsertParagraphAfter()

# Word keeps only a single "_GoBack" bookmark and relocates it to the most
# recently edited spot, so drop the old one before re-adding it on the new
# paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-fetch the freshly created (now-next) paragraph and populate it.
$newPara = $target.Next()

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:color w:val="FF0000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">When exiting full screen mode after watching a video, the arrows to go to the next slide reposition wrongly. They appear on top instead of in the middle. (on Firefox </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:color w:val="FF0000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>40.0.2</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
      <w:color w:val="FF0000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$null = $newPara.Range.InsertXML($xml)
